$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 49 (A49:G49) is the last existing entry for date "6.4.2020" / Mem Ctrl 2 / Add packages.
# Add a new row 50 with the same Task/Unit/Revision, extending the time-tracking table.

$ws.Range("A50").Value = $ws.Range("A49").Value
$ws.Range("A50").Style = $ws.Range("A49").Style

$ws.Range("B50").Value = 0.5
$ws.Range("B50").Style = $ws.Range("B49").Style

$ws.Range("C50").Value = 0.51388888888888895
$ws.Range("C50").Style = $ws.Range("C49").Style

$ws.Range("D50").Formula = "=C50-B50"
$ws.Range("D50").Style = $ws.Range("D49").Style

$ws.Range("E50").Value = $ws.Range("E48").Value
$ws.Range("E50").Style = $ws.Range("E49").Style

$ws.Range("F50").Value = $ws.Range("F49").Value
$ws.Range("F50").Style = $ws.Range("F49").Style

$ws.Range("G50").Value = $ws.Range("G49").Value
$ws.Range("G50").Style = $ws.Range("G49").Style

# Update the selection to reflect the new last-edited cell, as in the saved file.
$ws.Range("G50").Select()
